# Auto-generated edit script: updates leve-profit calculation cells
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 10002
$ws.Cells.Item(32, 10).Value = 10002
$ws.Cells.Item(32, 12).Value = 10002
$ws.Cells.Item(32, 14).Value = -10654

$ws.Cells.Item(86, 8).Value = 1000
$ws.Cells.Item(86, 10).Value = 1000
$ws.Cells.Item(86, 12).Value = 1000
$ws.Cells.Item(86, 14).Value = -3246

$ws.Cells.Item(89, 8).Value = 1000
$ws.Cells.Item(89, 10).Value = 1000
$ws.Cells.Item(89, 12).Value = 5000
$ws.Cells.Item(89, 14).Value = -16232

$ws.Cells.Item(129, 8).Value = 3031910.8
$ws.Cells.Item(129, 9).Value = 33333796
$ws.Cells.Item(129, 10).Value = 1722
$ws.Cells.Item(129, 11).Value = 100001388
$ws.Cells.Item(129, 12).Value = 5166
$ws.Cells.Item(129, 13).Value = -99996388
$ws.Cells.Item(129, 14).Value = -15166

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14743.143
$ws.Cells.Item(32, 9).Value = 14743.143
$ws.Cells.Item(32, 11).Value = 14743.143
$ws.Cells.Item(32, 13).Value = -14456.143

$ws.Cells.Item(61, 8).Value = 18498.5
$ws.Cells.Item(61, 9).Value = 17999
$ws.Cells.Item(61, 10).Value = 18665
$ws.Cells.Item(61, 11).Value = 17999
$ws.Cells.Item(61, 12).Value = 18665
$ws.Cells.Item(61, 13).Value = -17787
$ws.Cells.Item(61, 14).Value = -19089

$ws.Cells.Item(74, 8).Value = 6840
$ws.Cells.Item(74, 9).Value = 3540
$ws.Cells.Item(74, 11).Value = 3540
$ws.Cells.Item(74, 13).Value = -2666

$ws.Cells.Item(77, 8).Value = 6840
$ws.Cells.Item(77, 9).Value = 3540
$ws.Cells.Item(77, 11).Value = 17700
$ws.Cells.Item(77, 13).Value = -13332

$ws.Cells.Item(122, 8).Value = 1180
$ws.Cells.Item(122, 9).Value = 1250
$ws.Cells.Item(122, 10).Value = 900
$ws.Cells.Item(122, 11).Value = 3750
$ws.Cells.Item(122, 12).Value = 2700
$ws.Cells.Item(122, 13).Value = -1300
$ws.Cells.Item(122, 14).Value = -7600

$ws.Cells.Item(132, 8).Value = 7436.25
$ws.Cells.Item(132, 9).Value = 4581.6665
$ws.Cells.Item(132, 10).Value = 16000
$ws.Cells.Item(132, 11).Value = 13744.9995
$ws.Cells.Item(132, 12).Value = 48000
$ws.Cells.Item(132, 13).Value = -11214.9995
$ws.Cells.Item(132, 14).Value = -53060

$ws.Cells.Item(136, 8).Value = 18498.5
$ws.Cells.Item(136, 9).Value = 17999
$ws.Cells.Item(136, 10).Value = 18665
$ws.Cells.Item(136, 11).Value = 53997
$ws.Cells.Item(136, 12).Value = 55995
$ws.Cells.Item(136, 13).Value = -51447
$ws.Cells.Item(136, 14).Value = -61095

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 212.625
$ws.Cells.Item(22, 9).Value = 212.625
$ws.Cells.Item(22, 11).Value = 212.625
$ws.Cells.Item(22, 13).Value = -39.625

$ws.Cells.Item(36, 8).Value = 6444.25
$ws.Cells.Item(36, 9).Value = 6444.25
$ws.Cells.Item(36, 11).Value = 6444.25
$ws.Cells.Item(36, 13).Value = -5910.25

$ws.Cells.Item(94, 8).Value = 555.125
$ws.Cells.Item(94, 9).Value = 555.125
$ws.Cells.Item(94, 11).Value = 555.125
$ws.Cells.Item(94, 13).Value = -104.125

$ws.Cells.Item(103, 8).Value = 19862.285
$ws.Cells.Item(103, 10).Value = 19862.285
$ws.Cells.Item(103, 12).Value = 19862.285
$ws.Cells.Item(103, 14).Value = -22206.285

$ws.Cells.Item(134, 8).Value = 3712.4
$ws.Cells.Item(134, 9).Value = 1640.5
$ws.Cells.Item(134, 11).Value = 4921.5
$ws.Cells.Item(134, 13).Value = -2386.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 265.92307
$ws.Cells.Item(7, 9).Value = 284.75
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = 284.75
$ws.Cells.Item(7, 12).Value = 40
$ws.Cells.Item(7, 13).Value = -171.75
$ws.Cells.Item(7, 14).Value = -266

$ws.Cells.Item(58, 8).Value = 9681.5
$ws.Cells.Item(58, 9).Value = 5522.25
$ws.Cells.Item(58, 10).Value = 18000
$ws.Cells.Item(58, 11).Value = 5522.25
$ws.Cells.Item(58, 12).Value = 18000
$ws.Cells.Item(58, 13).Value = -5319.25
$ws.Cells.Item(58, 14).Value = -18406

$ws.Cells.Item(105, 8).Value = 2005
$ws.Cells.Item(105, 9).Value = 2005
$ws.Cells.Item(105, 11).Value = 2005
$ws.Cells.Item(105, 13).Value = -258

$ws.Cells.Item(107, 8).Value = 985.5
$ws.Cells.Item(107, 9).Value = 1418.3334
$ws.Cells.Item(107, 11).Value = 1418.3334
$ws.Cells.Item(107, 13).Value = 501.6666

$ws.Cells.Item(136, 8).Value = 9681.5
$ws.Cells.Item(136, 9).Value = 5522.25
$ws.Cells.Item(136, 10).Value = 18000
$ws.Cells.Item(136, 11).Value = 16566.75
$ws.Cells.Item(136, 12).Value = 54000
$ws.Cells.Item(136, 13).Value = -14016.75
$ws.Cells.Item(136, 14).Value = -59100

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 13).Value = ""

$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 13).Value = ""

$ws.Cells.Item(137, 8).Value = 2016
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 2016
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 6048
$ws.Cells.Item(137, 14).Value = -16248
$ws.Cells.Item(137, 13).Value = ""

$ws.Cells.Item(140, 8).Value = 4496.5
$ws.Cells.Item(140, 9).Value = 4496.5
$ws.Cells.Item(140, 11).Value = 13489.5
$ws.Cells.Item(140, 13).Value = -8309.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 2530
$ws.Cells.Item(41, 9).Value = 2530
$ws.Cells.Item(41, 11).Value = 2530
$ws.Cells.Item(41, 13).Value = -2175

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 5000
$ws.Cells.Item(17, 10).Value = 5000
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 14).Value = -5340

$ws.Cells.Item(22, 8).Value = 975
$ws.Cells.Item(22, 9).Value = 975
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 975
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -680
$ws.Cells.Item(22, 14).Value = ""

$ws.Cells.Item(27, 8).Value = 975
$ws.Cells.Item(27, 9).Value = 975
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 975
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -868
$ws.Cells.Item(27, 14).Value = ""

$ws.Cells.Item(30, 8).Value = 1261
$ws.Cells.Item(30, 9).Value = 613.2
$ws.Cells.Item(30, 10).Value = 4500
$ws.Cells.Item(30, 11).Value = 613.2
$ws.Cells.Item(30, 12).Value = 4500
$ws.Cells.Item(30, 13).Value = -505.2
$ws.Cells.Item(30, 14).Value = -4716

$ws.Cells.Item(46, 8).Value = 5600
$ws.Cells.Item(46, 9).Value = 5600
$ws.Cells.Item(46, 11).Value = 5600
$ws.Cells.Item(46, 13).Value = -5412

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 35814
$ws.Cells.Item(68, 10).Value = 35814
$ws.Cells.Item(68, 12).Value = 35814
$ws.Cells.Item(68, 14).Value = -37436

$ws.Cells.Item(71, 8).Value = 35814
$ws.Cells.Item(71, 10).Value = 35814
$ws.Cells.Item(71, 12).Value = 107442
$ws.Cells.Item(71, 14).Value = -115554

$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).Value = ""

$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).Value = ""

$ws.Cells.Item(107, 8).Value = 1389.7
$ws.Cells.Item(107, 10).Value = 1400.4
$ws.Cells.Item(107, 12).Value = 4201.200000000001
$ws.Cells.Item(107, 14).Value = -8041.200000000001

$ws.Cells.Item(132, 8).Value = 908
$ws.Cells.Item(132, 9).Value = 908
$ws.Cells.Item(132, 11).Value = 2724
$ws.Cells.Item(132, 13).Value = -194

$ws.Cells.Item(136, 8).Value = 7199.6
$ws.Cells.Item(136, 9).Value = 1332.6666
$ws.Cells.Item(136, 11).Value = 3997.9998
$ws.Cells.Item(136, 13).Value = -1447.9998
